$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "97.236.25"
$ws.Range("E2").Value = "  +1.64%  "

$ws.Range("D3").Value = "3.595.42"
$ws.Range("E3").Value = "  +0.29%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.87"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.73%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.77"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +16.43%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "653.20"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.44%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.426"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.73%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.02%  "

$ws.Range("E10").Value = "  +2.74%  "

$ws.Range("D11").Value = "3.593.63"
$ws.Range("E11").Value = "  +0.33%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "44.88"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.42%  "

$ws.Range("E13").Value = "  +0.92%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.48"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.00%  "

$ws.Range("D15").Value = "4.269.60"
$ws.Range("E15").Value = "  +0.43%  "

$ws.Range("D16").Value = "97.398.35"
$ws.Range("E16").Value = "  +1.97%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000261"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.27%  "

$ws.Range("D18").Value = "3.588.69"
$ws.Range("E18").Value = "  +0.35%  "

$ws.Range("E19").Value = "  +0.23%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.70"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.72%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "18.27"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.48%  "

$ws.Range("E22").Value = "  +6.64%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.48"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.36%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "517.48"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.32%  "

$ws.Range("E25").Value = "  +4.53%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.99"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.10%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "103.19"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +7.60%  "

$ws.Range("E28").Value = "  +3.37%  "

$ws.Range("E29").Value = "  +20.46%  "

$ws.Range("D30").Value = "3.793.42"
$ws.Range("E30").Value = "  +0.46%  "

$ws.Range("E31").Value = "  -2.19%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.98"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.09%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.22%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.191"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.95%  "

$ws.Range("E35").Value = "  +0.47%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "31.84"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.51%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.585"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.48%  "

$ws.Range("E38").Value = "  -1.42%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "616.29"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.98%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.76"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.35%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.155"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.12%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.93"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.37%  "

$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.462"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +44.22%  "

$ws.Range("E44").Value = "  -0.06%  "

$ws.Range("B45").Value = "ARBITRUM"
$ws.Range("C45").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.931"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.70%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.08"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.20%  "

$ws.Range("E47").Value = "  +6.67%  "

$ws.Range("E48").Value = "  +1.07%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "23.65"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.92%  "

$ws.Range("E50").Value = "  +5.21%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "33.08"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.36%  "
